$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: add a (collapsed) bookmark right at the very end of a
# paragraph's text (i.e. immediately before its paragraph mark). The
# runtime mis-handles Bookmarks.Add when given a range collapsed exactly
# at "paragraph.End - 1", so we work around it by temporarily inserting
# a marker character, anchoring the bookmark just before it, and then
# removing the marker again.
# -----------------------------------------------------------------------
function Add-BookmarkAtParagraphEnd($doc, $paraIndex, $name) {
    $para = $doc.Paragraphs.Item($paraIndex)
    $insPt = $doc.Range($para.Range.End - 1, $para.Range.End - 1)
    $insPt.InsertAfter("X")
    $para2 = $doc.Paragraphs.Item($paraIndex)
    $bookmarkPos = $para2.Range.End - 2
    $doc.Bookmarks.Add($name, $doc.Range($bookmarkPos, $bookmarkPos))
    $para3 = $doc.Paragraphs.Item($paraIndex)
    $xPos = $para3.Range.End - 1
    $doc.Range($xPos - 1, $xPos).Delete()
}

# ===========================================================================
# Change 1: first body paragraph - a few in-sentence word swaps / insertions
# around "sistem khusus ... menangani ... menggali ... Monitoring ... proses"
# ===========================================================================

# "sistem khusus yang menangani" -> "sistem khusus dalam menangani"
$r = $d.Content
$null = $r.Find.Execute("sistem khusus yang menangani", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordRange = $d.Range($r.Start + "sistem khusus ".Length, $r.End - " menangani".Length)
$wordRange.Text = "dalam"

# "menggali suatu informasi" -> "menggali dan mengelola suatu informasi"
$r2 = $d.Content
$null = $r2.Find.Execute("menggali suatu informasi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint2 = $d.Range($r2.Start + "menggali ".Length, $r2.Start + "menggali ".Length)
$insPoint2.InsertBefore("dan mengelola ")

# "Monitoring adalah suatu proses" -> "Monitoring berarti suatu proses"
$r3 = $d.Content
$null = $r3.Find.Execute("Monitoring adalah suatu proses", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$adalahRange = $d.Range($r3.Start + "Monitoring ".Length, $r3.End - " suatu proses".Length)
$adalahRange.Text = "berarti"

# ===========================================================================
# Change 2: "dapat diatas (W.H.O)" -> "dapat diatasi (W.H.O)"
# ===========================================================================
$r4 = $d.Content
$null = $r4.Find.Execute("dapat diatas (W.H.O)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$iPoint = $d.Range($r4.Start + "dapat diatas".Length, $r4.Start + "dapat diatas".Length)
$iPoint.InsertBefore("i")

# ===========================================================================
# Change 3: merge the empty paragraph right before "Matrikulasi ..." with
# that paragraph, expand it into the new, much longer STEI Tazkia /
# Matrikulasi discussion, and relocate the _GoBack bookmark to the very end
# of that paragraph.
# ===========================================================================

# Find the paragraph that starts with "Matrikulasi adalah" and the (empty)
# paragraph right before it, then merge them by deleting the empty
# paragraph's end-of-paragraph mark.
$matParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Matrikulasi adalah*") {
        $matParaIdx = $i
        break
    }
}
$emptyParaIdx = $matParaIdx - 1
$emptyPara = $d.Paragraphs.Item($emptyParaIdx)
$markRange = $d.Range($emptyPara.Range.End - 1, $emptyPara.Range.End)
$markRange.Delete()

# Insert the new introductory sentences before "Matrikulasi adalah ..."
$rM = $d.Content
$null = $rM.Find.Execute("Matrikulasi adalah", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPointM = $d.Range($rM.Start, $rM.Start)
$insPointM.InsertBefore("STEI Tazkia adalah suatu instansi pendidikan yang menerapkan program matrikulasi pada kegiatan pembelajaran dalam dua semester pertamanya. ")

# "Matrikulasi adalah kegiatan pembelajaran" -> "Matrikulasi merupakan kegiatan pembelajaran"
$rM2 = $d.Content
$null = $rM2.Find.Execute("Matrikulasi adalah kegiatan pembelajaran", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$adalahStart = $rM2.Start + "Matrikulasi ".Length
$adalahEnd = $adalahStart + "adalah".Length
$d.Range($adalahStart, $adalahEnd).Text = "merupakan"

# Append the large trailing discussion after "... diikuti [1]."
$rM3 = $d.Content
$null = $rM3.Find.Execute("diikuti [1].", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPointM = $d.Range($rM3.End, $rM3.End)
$trailing = " Program tersebut diterapkan dalam bentuk kegiatan boardingschool yang berjalan di Kampus Matrikulasi STEI Tazkia. Dalam tahap matrikulasi ini mahasiswa diwajibkan mengikuti berbagai program (sebagai syarat lulus tahap matrikulasi) didalamnya yaitu program Pembinaan, program Akademik dan program Bahasa (TLC / Tazkia Language Center). Ketiga program tersebut haruslah di monitor dengan baik oleh pihak manajemen matrikuklasi agar nantinya data bisa diolah dengan baik hingga dapat dijadikan suatu informasi yang mudah dibaca oleh pihak berkepentingan. Pada kenyataannya, seluruh kegiatan pada program tersebut belum ada suatu sistem yang menangani, dengan begitu data yang masuk hingga data yang sedang diolah menjadi informasi yang akan diterbitkan memiliki banyak kekurangan."
$endPointM.InsertBefore($trailing)

# Italicise "boardingschool " (with its trailing space)
$rItalic = $d.Content
$null = $rItalic.Find.Execute("boardingschool ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rItalic.Font.Italic = $true

# Move the _GoBack bookmark to the end of the (now much longer) paragraph.
$matParaIdx2 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "STEI Tazkia adalah*") {
        $matParaIdx2 = $i
        break
    }
}
Add-BookmarkAtParagraphEnd $d $matParaIdx2 "_GoBack"

foreach ($p in $d.Paragraphs) {
    Write-Host "==="
    Write-Host $p.Range.Text
}
